# Update "想去人数" (F column) values to match the latest scrape output.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 0
$ws1.Range("F4").Value = 0
$ws1.Range("F5").Value = 0
$ws1.Range("F8").Value = 148
$ws1.Range("F10").Value = 0

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 0
$ws4.Range("F5").Value = 0
$ws4.Range("F6").Value = 0
$ws4.Range("F7").Value = 0
$ws4.Range("F8").Value = 148
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 521
